$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header text change: single-layer theoretical latency unit ms -> us
$ws.Range("N3").Value = "单层理论延时(us)"

# 2. Update percentage ("占比(%)") values in column P for existing rows (4-17)
# because a new op row is being added and totals are recalculated.
$ws.Range("P4").Value = 3.32
$ws.Range("P5").Value = 2.07
$ws.Range("P6").Value = 0.92
$ws.Range("P7").Value = 6.43
$ws.Range("P8").Value = 2.85
$ws.Range("P9").Value = 1.42
$ws.Range("P10").Value = 1.53
$ws.Range("P11").Value = 12.23
$ws.Range("P12").Value = 6.11
$ws.Range("P13").Value = 6.11
$ws.Range("P14").Value = 3.06
$ws.Range("P15").Value = 7.35
$ws.Range("P16").Value = 3.67
$ws.Range("P17").Value = 7.35

# 3. Insert a new row at position 19. This pushes the old "combine" row
# (and everything below it) down by one, while the old "dispatch" row
# stays put at row 18 (its content will be replaced below to become the
# new "attn_all_reduce" row).
$ws.Rows.Item(19).Insert()

# Copy formatting (borders/number formats/etc.) from row 18 into the
# newly-inserted blank row 19 so it matches the rest of the data table.
$ws.Range("A18:P18").Copy()
$ws.Range("A19:P19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Row 18 becomes "attn_all_reduce" (repurposed from the old "dispatch" row)
$ws.Range("A18").Value = "attn_all_reduce"
$ws.Range("B18").Value = "transfer"
$ws.Range("C18").Value = 4096
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 7168
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 61
$ws.Range("H18").Value = "FP16"
$ws.Range("I18").Value = "FP16"
$ws.Range("J18").Value = "FP16"
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 690.827
$ws.Range("N18").Value = 690.827
$ws.Range("O18").Value = 42.14
$ws.Range("P18").Value = 9.24

# 5. Row 19 (the newly inserted row) becomes the "dispatch" row
$ws.Range("A19").Value = "dispatch"
$ws.Range("B19").Value = "transfer"
$ws.Range("C19").Value = 1024
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 7168
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = "INT8"
$ws.Range("I19").Value = "FP16"
$ws.Range("J19").Value = "FP16"
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 690.827
$ws.Range("N19").Value = 690.827
$ws.Range("O19").Value = 40.068
$ws.Range("P19").Value = 8.779999999999999

# 6. Row 20 ("combine", shifted down from old row 19): update D (n) and
# the recalculated percentage share.
$ws.Range("D20").Value = 0
$ws.Range("P20").Value = 17.56

# 7. Summary rows (shifted down by one row because of the insert above).
# 传输时间 (ms) value changes (row 26 now, was row 25).
$ws.Range("B26").Value = 2.763
# 总耗时 (ms) value changes (row 27 now, was row 26).
$ws.Range("B27").Value = 13.869

# 8. TTFT (ms) now at row 32 (was row 31) - updated value.
$ws.Range("B32").Value = 465.412

# 9. 吞吐量TPS now at row 33 (was row 32) - updated value.
$ws.Range("B33").Value = 8800.800999999999
